$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (Female, Estonia)
$ws.Range("D8").Value = 1750
$ws.Range("P8").Value = 358.8
$ws.Range("R8").Value = 25.8
$ws.Range("U8").Value = "358.8 (" + [char]0x00B1 + "104.4)"
$ws.Range("V8").Value = "25.8% (" + [char]0x00B1 + "8.8%)"
$ws.Range("X8").Value = 3987.1
$ws.Range("Z8").Value = "3987.1(" + [char]0x00B1 + "1160.1)"

# Row 13 (Female, Iceland)
$ws.Range("D13").Value = 255
$ws.Range("P13").Value = 12.4
$ws.Range("R13").Value = 5.1
$ws.Range("S13").Value = 9
$ws.Range("U13").Value = "12.4 (" + [char]0x00B1 + "22.7)"
$ws.Range("V13").Value = "5.1% (" + [char]0x00B1 + "9.0%)"
$ws.Range("X13").Value = 852.2
$ws.Range("Y13").Value = 1560.2
$ws.Range("Z13").Value = "852.2(" + [char]0x00B1 + "1560.2)"

# Row 39 (Male, Iceland)
$ws.Range("D39").Value = 148
$ws.Range("P39").Value = 8.6
$ws.Range("R39").Value = 6.2
$ws.Range("S39").Value = 7.7
$ws.Range("U39").Value = "8.6 (" + [char]0x00B1 + "10.9)"
$ws.Range("V39").Value = "6.2% (" + [char]0x00B1 + "7.7%)"
$ws.Range("X39").Value = 1187.8
$ws.Range("Y39").Value = 1505.6
$ws.Range("Z39").Value = "1187.8(" + [char]0x00B1 + "1505.6)"

# Row 60 (Total, Estonia)
$ws.Range("D60").Value = 2216
$ws.Range("P60").Value = 479.8
$ws.Range("S60").Value = 9.800000000000001
$ws.Range("U60").Value = "479.8 (" + [char]0x00B1 + "145.6)"
$ws.Range("V60").Value = "27.6% (" + [char]0x00B1 + "9.8%)"
$ws.Range("X60").Value = 4371.8
$ws.Range("Y60").Value = 1326.6
$ws.Range("Z60").Value = "4371.8(" + [char]0x00B1 + "1326.6)"

# Row 65 (Total, Iceland)
$ws.Range("D65").Value = 403
$ws.Range("P65").Value = 21
$ws.Range("R65").Value = 5.5
$ws.Range("S65").Value = 8
$ws.Range("U65").Value = "21.0 (" + [char]0x00B1 + "31.2)"
$ws.Range("V65").Value = "5.5% (" + [char]0x00B1 + "8.0%)"
$ws.Range("X65").Value = 963.7
$ws.Range("Y65").Value = 1431.9
$ws.Range("Z65").Value = "963.7(" + [char]0x00B1 + "1431.9)"
